# Generate Report for Handoff
# Adds a new row for file "fda37902-e119-45b6-a739-f64a6c3d6fde.md" to the
# Overview, zh-cn and de-de sheets of the localization-status workbook,
# mirroring the existing row for 4e19bfb8-0f26-4a77-9e9b-95a793e0e462.md.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/663aeb1f96697f780ce97928c0340895fe096531/e2e/"
$newFile = "fda37902-e119-45b6-a739-f64a6c3d6fde.md"
$dateFmt = "yyyy-mm-dd HH:mm:ss"
$hyperlinkColor = 15570276   # OLE (BGR) encoding of RGB FF6495ED, matches workbook's custom HyperLink style

# ---------------------------------------------------------------------------
# Sheet "Overview" (row 3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(3, 1).Value = $newFile
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = "2016-08-24 06:40:42"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseUrl + $newFile, [Type]::Missing, [Type]::Missing, "e2e\" + $newFile)
$wsOverview.Cells.Item(3, 2).Font.Name = "Calibri"
$wsOverview.Cells.Item(3, 2).Font.Size = 11
$wsOverview.Cells.Item(3, 2).Font.Underline = $true
$wsOverview.Cells.Item(3, 2).Font.Color = $hyperlinkColor

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (row 3)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(3, 1).Value = $newFile
$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = "'False"
$wsZhCn.Cells.Item(3, 7).Value = "fda37902-e119-45b6-a739-f64a6c3d6fde.0328c5ea439551c39cbeac87550f6e6a60e12c45.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 8).Value = "2016-08-24 06:40:37"
$wsZhCn.Cells.Item(3, 8).NumberFormat = $dateFmt
$wsZhCn.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3, 11).NumberFormat = $dateFmt
$wsZhCn.Cells.Item(3, 13).Value = "'True"
$wsZhCn.Cells.Item(3, 15).Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $baseUrl + $newFile, [Type]::Missing, [Type]::Missing, $newFile)
$wsZhCn.Cells.Item(3, 1).Font.Name = "Calibri"
$wsZhCn.Cells.Item(3, 1).Font.Size = 11
$wsZhCn.Cells.Item(3, 1).Font.Underline = $true
$wsZhCn.Cells.Item(3, 1).Font.Color = $hyperlinkColor

$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (row 3)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(3, 1).Value = $newFile
$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = "'False"
$wsDeDe.Cells.Item(3, 7).Value = "fda37902-e119-45b6-a739-f64a6c3d6fde.0328c5ea439551c39cbeac87550f6e6a60e12c45.de-de.xlf"
$wsDeDe.Cells.Item(3, 8).Value = "2016-08-24 06:40:42"
$wsDeDe.Cells.Item(3, 8).NumberFormat = $dateFmt
$wsDeDe.Cells.Item(3, 11).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3, 11).NumberFormat = $dateFmt
$wsDeDe.Cells.Item(3, 13).Value = "'True"
$wsDeDe.Cells.Item(3, 15).Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $baseUrl + $newFile, [Type]::Missing, [Type]::Missing, $newFile)
$wsDeDe.Cells.Item(3, 1).Font.Name = "Calibri"
$wsDeDe.Cells.Item(3, 1).Font.Size = 11
$wsDeDe.Cells.Item(3, 1).Font.Underline = $true
$wsDeDe.Cells.Item(3, 1).Font.Color = $hyperlinkColor

$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
